$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Ticketart attribute definition: "select" -> "dropdown"
$ws.Range("E3").Value = "Ticketart:dropdown(Bus,Zug,U-Bahn);Häufigkeit:dropdown(Täglich,Wöchentlich,Selten)"

# Move the active selection to E3 (single cell), matching the saved view state
$ws.Activate()
$ws.Range("E3").Select()
